# Updates "想去人数" (want-to-go count, column F) and, for the two rows
# whose ticket became sold out, the "最低票价" (lowest price, column G)
# across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 393
$ws.Range("F6").Value = 929
$ws.Range("F7").Value = 4184
$ws.Range("F8").Value = 339
$ws.Range("F9").Value = 218
$ws.Range("F10").Value = 832
$ws.Range("F12").Value = 60
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 524
$ws.Range("F17").Value = 1503
$ws.Range("F18").Value = 1392
$ws.Range("F19").Value = 587
$ws.Range("F20").Value = 295
$ws.Range("F21").Value = 162
$ws.Range("F22").Value = 208
$ws.Range("F23").Value = 420
$ws.Range("F24").Value = 85
$ws.Range("F25").Value = 1029
$ws.Range("F26").Value = 16
$ws.Range("F27").Value = 520
$ws.Range("G27").Value = "已售罄"
$ws.Range("F28").Value = 844
$ws.Range("F29").Value = 103
$ws.Range("F30").Value = 66
$ws.Range("F31").Value = 143
$ws.Range("F35").Value = 226
$ws.Range("F36").Value = 229
$ws.Range("F37").Value = 486
$ws.Range("F38").Value = 52

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 116
$ws.Range("F6").Value = 86

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 226

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - combination of the other three sheets
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 226
$ws.Range("F6").Value = 393
$ws.Range("F7").Value = 929
$ws.Range("F9").Value = 4184
$ws.Range("F10").Value = 339
$ws.Range("F11").Value = 218
$ws.Range("F13").Value = 116
$ws.Range("F14").Value = 832
$ws.Range("F17").Value = 86
$ws.Range("F18").Value = 60
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 524
$ws.Range("F24").Value = 1503
$ws.Range("F25").Value = 1392
$ws.Range("F26").Value = 587
$ws.Range("F27").Value = 295
$ws.Range("F28").Value = 162
$ws.Range("F29").Value = 208
$ws.Range("F31").Value = 420
$ws.Range("F32").Value = 85
$ws.Range("F33").Value = 1029
$ws.Range("F34").Value = 16
$ws.Range("F35").Value = 520
$ws.Range("G35").Value = "已售罄"
$ws.Range("F36").Value = 844
$ws.Range("F37").Value = 103
$ws.Range("F38").Value = 66
$ws.Range("F39").Value = 143
$ws.Range("F43").Value = 226
$ws.Range("F44").Value = 229
$ws.Range("F45").Value = 486
$ws.Range("F46").Value = 52
